$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices, percentages) stay as text,
# matching the original inline-string cell content exactly.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.963.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.745.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.85%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.18%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.744.49"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.85%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.76%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.49"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.55"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.377.53"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.92%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.747.55"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.14%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.001.96"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.05"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.36%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.112"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.80"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.61"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.703"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.18"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.62%  "

# Row 25
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000137"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -9.01%  "

# Row 26
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.55%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.71%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.61%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.895.54"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.85%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.38"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.24"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.91"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.49%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.12"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.94%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.701.22"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.23%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -11.29%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.27%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.992"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.98%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.00%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.308"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.19%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.60"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.76%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.92"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.57%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.36"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "395.18"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.08%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.26"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.61%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.59"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.53%  "
